$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B51").Value = "U22"
$ws.Range("A51").Value = "DRV8833"
$ws.Range("D51").Value = "C50506"
$ws.Range("C51").Value = "Package_SO:HTSSOP-16-1EP_4.4x5mm_P0.65mm_EP3.4x5mm_Mask2.46x2.31mm_ThermalVias"
